$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 6: CONCAT formula spanning SPAM literals and the A1:B2 range
$ws.Range("A6").Formula = '=CONCAT("SPAM", " ", A1:B2, "SPAM", " ")'

# Move the active selection to A7 (matches the post-edit cursor position)
$ws.Range("A7").Select() | Out-Null
